# Update crypto price (D) and volume/change (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.721.81"
$ws.Range("D2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.798.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.22%  "
# Row 4
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "
# Row 6
$ws.Range("E6").Value = "  -0.06%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4470"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.70%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3681"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.40%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07345"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.69%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8611"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.20%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.800.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.30%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.631"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.23%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.13%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07070"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.14%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.273"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.15%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008680"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.69%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.739.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.69%  "
# Row 22
$ws.Range("E22").Value = "  +0.81%  "
# Row 23
$ws.Range("E23").Value = "  -0.30%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.018.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.75%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.983"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.166"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.23%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.198"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.91%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08779"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.19%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7400"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.158"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.88%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.455"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.914"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.64%  "
# Row 36
$ws.Range("E36").Value = "  -0.06%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.083"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.27%  "
# Row 38
$ws.Range("E38").Value = "  -0.34%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05194"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.90%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5277"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.96%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.832"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.45%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.969"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.57%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1681"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5075"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.89%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.424"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.93%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.22%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.960"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.70%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.07%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.676"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.25%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9995"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.08%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06294"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
